# Update "想去人数" (interest count) figures in column F across the
# 展览 (Exhibition), 演出 (Performance) and 全部类型 (All types) sheets,
# matching the refreshed data snapshot for the gh-pages output.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 306
$ws1.Range("F4").Value = 2943
$ws1.Range("F5").Value = 75
$ws1.Range("F7").Value = 2306
$ws1.Range("F8").Value = 1624
$ws1.Range("F10").Value = 845
$ws1.Range("F11").Value = 106
$ws1.Range("F12").Value = 2644
$ws1.Range("F14").Value = 1493
$ws1.Range("F15").Value = 6968
$ws1.Range("F17").Value = 7140
$ws1.Range("F19").Value = 4838
$ws1.Range("F20").Value = 3089
$ws1.Range("F21").Value = 3459
$ws1.Range("F24").Value = 1841
$ws1.Range("F28").Value = 14
$ws1.Range("F29").Value = 173
$ws1.Range("F31").Value = 1978
$ws1.Range("F32").Value = 1119
$ws1.Range("F33").Value = 2586
$ws1.Range("F34").Value = 5
$ws1.Range("F36").Value = 160
$ws1.Range("F37").Value = 374
$ws1.Range("F38").Value = 1028
$ws1.Range("F39").Value = 205
$ws1.Range("F40").Value = 466

# --- Sheet: 演出 (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F9").Value = 31
$ws2.Range("F14").Value = 91

# --- Sheet: 全部类型 (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 306
$ws4.Range("F6").Value = 2943
$ws4.Range("F7").Value = 75
$ws4.Range("F8").Value = 2306
$ws4.Range("F9").Value = 1624
$ws4.Range("F11").Value = 845
$ws4.Range("F12").Value = 106
$ws4.Range("F14").Value = 2644
$ws4.Range("F15").Value = 1493
$ws4.Range("F17").Value = 31
$ws4.Range("F19").Value = 6968
$ws4.Range("F21").Value = 7140
$ws4.Range("F23").Value = 4839
$ws4.Range("F24").Value = 3089
$ws4.Range("F25").Value = 3459
$ws4.Range("F29").Value = 91
$ws4.Range("F30").Value = 1841
$ws4.Range("F35").Value = 14
$ws4.Range("F36").Value = 173
$ws4.Range("F38").Value = 1979
$ws4.Range("F39").Value = 1119
$ws4.Range("F41").Value = 2586
$ws4.Range("F43").Value = 160
$ws4.Range("F45").Value = 374
$ws4.Range("F46").Value = 1028
$ws4.Range("F47").Value = 205
$ws4.Range("F48").Value = 466
